$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C2:C27) from 2023-12-26 (45286)
# to 2023-12-30 (45290) for every data row.
$ws.Range("C2:C27").Value = 45290
